# Update "想去人数" (F column) values across sheets to reflect the latest
# scraped counts, as produced by the gh-pages data generation run at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12833
$ws1.Range("F3").Value = 625
$ws1.Range("F6").Value = 321
$ws1.Range("F7").Value = 402
$ws1.Range("F9").Value = 12829
$ws1.Range("F11").Value = 20
$ws1.Range("F12").Value = 5232
$ws1.Range("F14").Value = 18
$ws1.Range("F18").Value = 36
$ws1.Range("F21").Value = 2851
$ws1.Range("F22").Value = 6171
$ws1.Range("F23").Value = 1155
$ws1.Range("F24").Value = 3623
$ws1.Range("F25").Value = 221

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 24

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12833
$ws4.Range("F3").Value = 625
$ws4.Range("F6").Value = 321
$ws4.Range("F7").Value = 24
$ws4.Range("F8").Value = 402
$ws4.Range("F10").Value = 12829
$ws4.Range("F12").Value = 20
$ws4.Range("F13").Value = 5232
$ws4.Range("F15").Value = 18
$ws4.Range("F19").Value = 36
$ws4.Range("F22").Value = 2851
$ws4.Range("F24").Value = 6171
$ws4.Range("F25").Value = 1155
$ws4.Range("F26").Value = 3623
$ws4.Range("F27").Value = 221
